$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 7 ("Experimental"): B7 was blank, now holds the literal text "false"
# (this must remain a text string, not a Boolean, to match the generated
# report's shared-string cell). A bare Value = "false" gets auto-typed to
# a Boolean by Excel, so a leading apostrophe is used to force literal-text
# entry. That marks the cell as "quote prefixed", which would leave a
# stray quote-prefix format on B7, so the original cell format (shared
# with every other value cell in this column) is restored afterwards via
# copy/paste-special of a neighbouring cell's formatting only.
$b7 = $ws.Cells.Item(7, 2)
$b7.ClearFormats()
$b7.Value = "'false"
$ws.Cells.Item(6, 2).Copy()
$b7.PasteSpecial(-4122)

# Row 8 ("Date"): refresh the generation timestamp
$ws.Cells.Item(8, 2).Value = "2025-11-30T13:08:37+00:00"
